$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 28; $r -le 176; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
